$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: "Year" header (C1) gets a numeric-style number format ---
$ws.Range("C1").NumberFormat = "#,##0"

# --- Row 2: Year becomes a real number instead of text "2026" ---
$ws.Range("C2").Value = 2026
$ws.Range("C2").NumberFormat = "#,##0"

# --- Row 4: was a placeholder (blank/quote-prefixed) row; fill with the new Circulars record ---
$ws.Range("A4").Value = "SEBI"
$ws.Range("B4").Value = "Circulars"
$ws.Range("C4").Value = 2026
$ws.Range("C4").NumberFormat = "#,##0"
$ws.Range("D4").Value = "January"
# Force the IssueDate to stay literal text (not auto-parsed into a date serial),
# then restore the plain body formatting that a normal text entry would carry.
$ws.Range("E4").Value = "'2026-01-02"
$ws.Range("D4").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("F4").Value = "Specification of the consequential requirements with respect to Amendment of Securities and Exchange Board of India (Merchant Bankers) Regulations, 1992"
$ws.Range("G4").Value = "https://www.sebi.gov.in/sebi_data/attachdocs/jan-2026/1767358255887.pdf"
$ws.Range("H4").Value = "1767358255887.pdf"
$ws.Range("I4").Value = "/Users/admin/Downloads/Tejomaya_pdfs/Akshayam Data/SEBI/Circulars/2026/January/1767358255887.pdf"

# --- Rows 5 & 6: drop the stray quote-prefixed empty placeholders and leave true blank cells ---
$ws.Range("A5:I6").Value = "x"
$ws.Range("C5:C6").NumberFormat = "#,##0"
$ws.Range("A5:I6").ClearContents()
